$wb = $excel.ActiveWorkbook

function Set-TopBottomBorder($rng) {
  # Turn off left/right, turn on thin top+bottom (matches styles xf borderId=4)
  $rng.Borders.Item(7).LineStyle = 0
  $rng.Borders.Item(10).LineStyle = 0
  $rng.Borders.Item(8).LineStyle = 1
  $rng.Borders.Item(8).Weight = 2
  $rng.Borders.Item(9).LineStyle = 1
  $rng.Borders.Item(9).Weight = 2
}

function Set-TopBottomRightBorder($rng) {
  # Turn off left, turn on thin top+bottom+right (matches styles xf borderId=5)
  $rng.Borders.Item(7).LineStyle = 0
  $rng.Borders.Item(8).LineStyle = 1
  $rng.Borders.Item(8).Weight = 2
  $rng.Borders.Item(9).LineStyle = 1
  $rng.Borders.Item(9).Weight = 2
  $rng.Borders.Item(10).LineStyle = 1
  $rng.Borders.Item(10).Weight = 2
}

# ----- Sheet 1: quality_comparison -----
$ws1 = $wb.Worksheets.Item("quality_comparison")

Set-TopBottomBorder($ws1.Range("C1"))
Set-TopBottomRightBorder($ws1.Range("D1"))

$ws1.Range("C2").Value = "approach"

# ----- Sheet 2: computational_comparison -----
$ws2 = $wb.Worksheets.Item("computational_comparison")

Set-TopBottomBorder($ws2.Range("C1"))
Set-TopBottomRightBorder($ws2.Range("D1"))
Set-TopBottomBorder($ws2.Range("F1"))
Set-TopBottomRightBorder($ws2.Range("G1"))

$ws2.Range("C2").Value = "approach"
$ws2.Range("F2").Value = "approach"

$ws2.Range("G5").ClearContents()
